$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = 44434
$ws.Range("M12").Value = 100
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 20000
$ws.Range("S12").Value = 2000
$ws.Range("D13").Value = 44340
$ws.Range("M13").Value = 45
$ws.Range("P13").Value = 20556
$ws.Range("S13").Value = 2056
$ws.Range("D14").Value = 44280
$ws.Range("M14").Value = 80
$ws.Range("P14").Value = 20500
$ws.Range("S14").Value = 2050
$ws.Range("D15").Value = 44362
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 21000
$ws.Range("P15").Value = 20600
$ws.Range("S15").Value = 2060
$ws.Range("D16").Value = 44242
$ws.Range("M16").Value = 55
$ws.Range("N16").Value = 25000
$ws.Range("O16").Value = 25000
$ws.Range("P16").Value = 25000
$ws.Range("S16").Value = 2500
$ws.Range("D17").Value = 44363
$ws.Range("N17").Value = 21000
$ws.Range("P17").Value = 21000
$ws.Range("S17").Value = 2100
$ws.Range("D18").Value = 44349
$ws.Range("M18").Value = 80
$ws.Range("N18").Value = 20000
$ws.Range("P18").Value = 20500
$ws.Range("S18").Value = 2050
$ws.Range("D19").Value = 44385
$ws.Range("M19").Value = 25
$ws.Range("N19").Value = 21000
$ws.Range("O19").Value = 21000
$ws.Range("P19").Value = 21000
$ws.Range("S19").Value = 2100
$ws.Range("D20").Value = 44258
$ws.Range("M20").Value = 15
$ws.Range("N20").Value = 22000
$ws.Range("O20").Value = 22000
$ws.Range("P20").Value = 22000
$ws.Range("S20").Value = 2200
$ws.Range("D21").Value = 44427
$ws.Range("N21").Value = 20000
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 20000
$ws.Range("Q21").Value = "$/bandeja 10 kilos"
$ws.Range("S21").Value = 2000
$ws.Range("T21").Value = 10
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 80
$ws.Range("N22").Value = 2500
$ws.Range("O22").Value = 2500
$ws.Range("P22").Value = 2500
$ws.Range("S22").Value = 2500
$ws.Range("D23").Value = 44413
$ws.Range("L23").Value = "Segunda"
$ws.Range("M23").Value = 55
$ws.Range("N23").Value = 2000
$ws.Range("O23").Value = 2000
$ws.Range("P23").Value = 2000
$ws.Range("Q23").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("T23").Value = 1
$ws.Range("D24").Value = 44421
$ws.Range("M24").Value = 30
$ws.Range("N24").Value = 20000
$ws.Range("O24").Value = 20000
$ws.Range("P24").Value = 20000
$ws.Range("S24").Value = 2000
$ws.Range("D25").Value = 44377
$ws.Range("M25").Value = 100
$ws.Range("D26").Value = 44321
$ws.Range("M26").Value = 80
$ws.Range("N26").Value = 21000
$ws.Range("O26").Value = 21000
$ws.Range("P26").Value = 21000
$ws.Range("Q26").Value = "$/bandeja 10 kilos"
$ws.Range("S26").Value = 2100
$ws.Range("T26").Value = 10
$ws.Range("D27").Value = 44412
$ws.Range("M27").Value = 25
$ws.Range("N27").Value = 2500
$ws.Range("O27").Value = 2500
$ws.Range("P27").Value = 2500
$ws.Range("Q27").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R27").Value = "Provincia del Elquí"
$ws.Range("S27").Value = 2500
$ws.Range("T27").Value = 1
$ws.Range("D28").Value = 44336
$ws.Range("M28").Value = 40
$ws.Range("R28").Value = "Provincia de Limarí"
$ws.Range("D29").Value = 44405
$ws.Range("M29").Value = 70
$ws.Range("N29").Value = 20000
$ws.Range("O29").Value = 20000
$ws.Range("P29").Value = 20000
$ws.Range("S29").Value = 2000
$ws.Range("D30").Value = 44435
$ws.Range("M30").Value = 390
$ws.Range("O30").Value = 21000
$ws.Range("P30").Value = 20128
$ws.Range("S30").Value = 2013
$ws.Range("D31").Value = 44431
$ws.Range("M31").Value = 130
$ws.Range("O31").Value = 21000
$ws.Range("P31").Value = 20385
$ws.Range("S31").Value = 2038
$ws.Range("D32").Value = 44251
$ws.Range("M32").Value = 55
$ws.Range("N32").Value = 24000
$ws.Range("O32").Value = 24000
$ws.Range("P32").Value = 24000
$ws.Range("S32").Value = 2400
$ws.Range("D33").Value = 44417
$ws.Range("M33").Value = 150
$ws.Range("N33").Value = 20000
$ws.Range("O33").Value = 20000
$ws.Range("P33").Value = 20000
$ws.Range("S33").Value = 2000
$ws.Range("D34").Value = 44419
$ws.Range("M34").Value = 100
$ws.Range("N34").Value = 20000
$ws.Range("O34").Value = 20000
$ws.Range("P34").Value = 20000
$ws.Range("Q34").Value = "$/bandeja 10 kilos"
$ws.Range("S34").Value = 2000
$ws.Range("T34").Value = 10
$ws.Range("D35").Value = 44307
$ws.Range("M35").Value = 50
$ws.Range("Q35").Value = "$/bandeja 10 kilos"
$ws.Range("S35").Value = 2100
$ws.Range("T35").Value = 10
$ws.Range("D36").Value = 44265
$ws.Range("M36").Value = 40
$ws.Range("N36").Value = 21000
$ws.Range("O36").Value = 21000
$ws.Range("P36").Value = 21000
$ws.Range("Q36").Value = "$/bandeja 10 kilos"
$ws.Range("S36").Value = 2100
$ws.Range("D37").Value = 44333
$ws.Range("M37").Value = 30
$ws.Range("N37").Value = 35000
$ws.Range("O37").Value = 35000
$ws.Range("P37").Value = 35000
$ws.Range("Q37").Value = "$/caja 15 kilos granel"
$ws.Range("S37").Value = 2333
$ws.Range("T37").Value = 15
$ws.Range("D38").Value = 44301
$ws.Range("M38").Value = 55
$ws.Range("Q38").Value = "$/caja 15 kilos granel"
$ws.Range("S38").Value = 1400
$ws.Range("T38").Value = 15
$ws.Range("D39").Value = 44433
$ws.Range("M39").Value = 80
$ws.Range("N39").Value = 20000
$ws.Range("O39").Value = 20000
$ws.Range("P39").Value = 20000
$ws.Range("Q39").Value = "$/bandeja 10 kilos"
$ws.Range("S39").Value = 2000
$ws.Range("T39").Value = 10
$ws.Range("D40").Value = 44370
$ws.Range("M40").Value = 20
$ws.Range("Q40").Value = "$/caja 10 kilos"

# New rows appended at the end
$row = 41
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44326
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100108
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value = 100108004
$ws.Cells.Item($row, 10).Value = "Papaya"
$ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 55
$ws.Cells.Item($row, 14).Value = 21000
$ws.Cells.Item($row, 15).Value = 21000
$ws.Cells.Item($row, 16).Value = 21000
$ws.Cells.Item($row, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 19).Value = 2100
$ws.Cells.Item($row, 20).Value = 10
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 42
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44382
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100108
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value = 100108004
$ws.Cells.Item($row, 10).Value = "Papaya"
$ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 45
$ws.Cells.Item($row, 14).Value = 21000
$ws.Cells.Item($row, 15).Value = 21000
$ws.Cells.Item($row, 16).Value = 21000
$ws.Cells.Item($row, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 19).Value = 2100
$ws.Cells.Item($row, 20).Value = 10
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 43
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44334
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100108
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value = 100108004
$ws.Cells.Item($row, 10).Value = "Papaya"
$ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 20
$ws.Cells.Item($row, 14).Value = 35000
$ws.Cells.Item($row, 15).Value = 35000
$ws.Cells.Item($row, 16).Value = 35000
$ws.Cells.Item($row, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item($row, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 19).Value = 2333
$ws.Cells.Item($row, 20).Value = 15
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 44
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44418
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100108
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value = 100108004
$ws.Cells.Item($row, 10).Value = "Papaya"
$ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 50
$ws.Cells.Item($row, 14).Value = 20000
$ws.Cells.Item($row, 15).Value = 20000
$ws.Cells.Item($row, 16).Value = 20000
$ws.Cells.Item($row, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 19).Value = 2000
$ws.Cells.Item($row, 20).Value = 10
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 45
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44432
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100108
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value = 100108004
$ws.Cells.Item($row, 10).Value = "Papaya"
$ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 80
$ws.Cells.Item($row, 14).Value = 20000
$ws.Cells.Item($row, 15).Value = 20000
$ws.Cells.Item($row, 16).Value = 20000
$ws.Cells.Item($row, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 19).Value = 2000
$ws.Cells.Item($row, 20).Value = 10
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
